$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price values stored as TEXT strings.
# A leading apostrophe forces Excel to keep them as text instead of
# auto-converting to a number; the Style reset afterwards drops the
# "quote prefix" formatting Excel applies so the cell keeps its original
# (default/unstyled) look.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "241.40"
Set-TextValue $ws.Range("D3") "21.76"
Set-TextValue $ws.Range("D4") "5.356"
Set-TextValue $ws.Range("D5") "0.05674"
Set-TextValue $ws.Range("D6") "3.407"
Set-TextValue $ws.Range("D7") "6.283"
Set-TextValue $ws.Range("D8") "0.8072"
Set-TextValue $ws.Range("D9") "0.8643"
Set-TextValue $ws.Range("D10") "0.1423"
Set-TextValue $ws.Range("D11") "0.07272"
Set-TextValue $ws.Range("D12") "0.03053"
Set-TextValue $ws.Range("D13") "0.03158"
Set-TextValue $ws.Range("D14") "0.09348"
Set-TextValue $ws.Range("D15") "3.906"
Set-TextValue $ws.Range("D16") "0.001584"
Set-TextValue $ws.Range("D17") "0.04811"
Set-TextValue $ws.Range("D18") "0.0005825"
Set-TextValue $ws.Range("D19") "0.006314"
Set-TextValue $ws.Range("D20") "0.0009971"
Set-TextValue $ws.Range("D21") "0.004053"
Set-TextValue $ws.Range("D22") "0.0001501"
Set-TextValue $ws.Range("D23") "3.738"
Set-TextValue $ws.Range("D24") "2.151"
Set-TextValue $ws.Range("D27") "0.0004003"
Set-TextValue $ws.Range("D40") "0.03788"
Set-TextValue $ws.Range("D41") "0.006687"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
Set-TextValue $ws.Range("D42") "0.1045"
Set-TextValue $ws.Range("D43") "0.002831"
Set-TextValue $ws.Range("D44") "0.006831"
Set-TextValue $ws.Range("D45") "0.00005608"
Set-TextValue $ws.Range("D47") "0.5805"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
Set-TextValue $ws.Range("D48") "0.1423"
Set-TextValue $ws.Range("D49") "0.00002102"
Set-TextValue $ws.Range("D50") "0.01011"
